$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New runs added to the results table (column A/B) and the raw softmax
# arrays tucked away in column F, mirroring the existing layout used for
# the prior runs (Normal, NM+CutMix(custom), NM+CutMix(rect), NM+Mixup).
# Written in this interleaved order so the shared-string table is built up
# in the same sequence as the source workbook.

$ws.Range("A6").Value = "NM(64)"
$ws.Range("F23").Value = "[0.12058988213539124, 0.03436800464987755, 0.10352993756532669, 0.02030286379158497, 0.03217571601271629, 0.02710207737982273, 0.18644094467163086, 0.3134118318557739, 0.25580230355262756, 0.13114772737026215, 0.1224871288985014]"

$ws.Range("A7").Value = "NM+4CM(rect)(64)"
$ws.Range("F24").Value = "[0.11653576046228409, 0.03403253108263016, 0.10798562318086624, 0.021862655878067017, 0.03326644003391266, 0.027488116174936295, 0.20256248116493225, 0.238773375749588, 0.23526504635810852, 0.1368890106678009, 0.11546610407531262]"

$ws.Range("A8").Value = "NM+4CM(custom)(64)"
$ws.Range("F25").Value = "[0.11634423583745956, 0.03415555879473686, 0.10137994587421417, 0.0216276366263628, 0.032765522599220276, 0.027110232040286064, 0.18748654425144196, 0.23630554974079132, 0.24979417026042938, 0.13126122951507568, 0.1138230625540018]"

$ws.Range("A9").Value = "NM+4MixUp(64)"
$ws.Range("F26").Value = "[0.11625339835882187, 0.033972881734371185, 0.10771539062261581, 0.024174481630325317, 0.03213963657617569, 0.027086645364761353, 0.17631056904792786, 0.3103159964084625, 0.2537173330783844, 0.12410194426774979, 0.12057882770895959]"

$ws.Range("F28").Value = "[0.11612879484891891, 0.03450104221701622, 0.1083734780550003, 0.027244193479418755, 0.032124631106853485, 0.026846332475543022, 0.19984659552574158, 0.25128018856048584, 0.22460828721523285, 0.13065199553966522, 0.11516055390238762]"
$ws.Range("A11").Value = "NM+9CM(custom)(64)"

# Numeric results (last element of each array above) in column B, matching
# the pattern already used for the first four runs.
$ws.Range("B6").Value = 0.122487128898501
$ws.Range("B7").Value = 0.11546610407531201
$ws.Range("B8").Value = 0.113823062554001
$ws.Range("B9").Value = 0.120578827708959
$ws.Range("B11").Value = 0.115160553902387

# Column A grew a couple more characters (longest label is now
# "NM+4CM(custom)(64)"), so re-fit it like the rest of the sheet.
$ws.Columns.Item(1).AutoFit()

# Leave the selection where the author's cursor ended up.
$null = $ws.Range("B11").Select()
